$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.273.83'
$ws.Range("E2").Value = '  +2.15%  '

$ws.Range("D3").Value = '3.399.78'
$ws.Range("E3").Value = '  +2.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.92'
$ws.Range("E5").Value = '  +1.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.24'
$ws.Range("E6").Value = '  +3.62%  '

$ws.Range("E7").Value = '  +1.82%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.201'
$ws.Range("E9").Value = '  +10.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.594'
$ws.Range("E10").Value = '  +2.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.58'
$ws.Range("E11").Value = '  +3.33%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000286'
$ws.Range("E12").Value = '  +5.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '684.32'
$ws.Range("E13").Value = '  -1.72%  '

$ws.Range("E14").Value = '  +4.00%  '

$ws.Range("D15").Value = '3.947.11'
$ws.Range("E15").Value = '  +2.21%  '

$ws.Range("D16").Value = '69.368.11'
$ws.Range("E16").Value = '  +2.34%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.408.56'
$ws.Range("E17").Value = '  +2.28%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.121'
$ws.Range("E18").Value = '  +1.62%  '

$ws.Range("E19").Value = '  +1.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.33'
$ws.Range("E20").Value = '  +2.11%  '

$ws.Range("E21").Value = '  +2.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.36'
$ws.Range("E22").Value = '  +2.75%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.37'
$ws.Range("E23").Value = '  -0.73%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '103.25'
$ws.Range("E24").Value = '  +1.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.95'
$ws.Range("E25").Value = '  +1.34%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.74'
$ws.Range("E26").Value = '  +2.09%  '

$ws.Range("E27").Value = '  +3.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.94'
$ws.Range("E28").Value = '  +3.16%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.81'
$ws.Range("E29").Value = '  +3.40%  '

$ws.Range("E30").Value = '  -0.25%  '

$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.17'
$ws.Range("E31").Value = '  +1.70%  '

$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '558.95'
$ws.Range("E32").Value = '  -2.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.65'
$ws.Range("E33").Value = '  +11.40%  '

$ws.Range("E34").Value = '  +1.66%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.69'
$ws.Range("E35").Value = '  +3.82%  '

$ws.Range("E36").Value = '  -0.03%  '

$ws.Range("D37").Value = '3.667.28'
$ws.Range("E37").Value = '  -1.38%  '

$ws.Range("E38").Value = '  +6.64%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.07'
$ws.Range("E39").Value = '  +1.51%  '

$ws.Range("D40").Value = '0.0₃0725'
$ws.Range("E40").Value = '  +8.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.26'
$ws.Range("E41").Value = '  +4.13%  '

$ws.Range("E42").Value = '  +2.94%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.339'
$ws.Range("E43").Value = '  +1.60%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0428'
$ws.Range("E44").Value = '  +5.56%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.34'
$ws.Range("E45").Value = '  +0.70%  '

$ws.Range("E46").Value = '  +2.04%  '

$ws.Range("E47").Value = '  +1.22%  '

$ws.Range("E48").Value = '  +5.34%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.95'
$ws.Range("E50").Value = '  +0.68%  '

$ws.Range("E51").Value = '  +2.96%  '
